$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain cell assignments for values that Excel's General format
# will not mis-parse as numbers/dates (URLs, names, multi-dot
# price strings, percent strings with leading/trailing spaces).
$ws.Range('D2').Value = '25.645.03'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '1.689.35'
$ws.Range('E3').Value = '  -2.83%  '
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('E7').Value = '  -5.59%  '
$ws.Range('E8').Value = '  -3.23%  '
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').Value = '1.716.65'
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('E12').Value = '  -1.19%  '
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('E15').Value = '  -3.70%  '
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('D18').Value = '25.657.20'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E21').Value = '  +5.27%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '1.921.31'
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('E25').Value = '  -3.56%  '
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('E27').Value = '  -7.63%  '
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('E31').Value = '  -3.10%  '
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('E33').Value = '  +0.94%  '
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('E35').Value = '  -1.76%  '
$ws.Range('E36').Value = '  -4.06%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').Value = '  -1.83%  '
$ws.Range('E39').Value = '  +14.55%  '
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('E42').Value = '  +2.82%  '
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('E45').Value = '  +2.65%  '
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('E49').Value = '  -1.28%  '
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('E51').Value = '  -0.51%  '

# The cells below hold plain decimal numbers (e.g. '1.005',
# '30.00') that must stay stored as literal TEXT, exactly as
# they were originally (matching the source feed's inline
# strings) instead of being auto-coerced to numeric values
# (which would silently drop trailing zeros / change type).
# A scratch cell is formatted as Text, written once per value,
# then its value-only (not formatting) is pasted onto the real
# target cell so the target keeps its original (General) style.
$scratch = $ws.Range('ZZ1')
$scratch.Value = "'1.005"
$scratch.Copy()
$ws.Range('D4').PasteSpecial(-4163)
$scratch.Value = "'240.90"
$scratch.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$scratch.Value = "'0.4894"
$scratch.Copy()
$ws.Range('D7').PasteSpecial(-4163)
$scratch.Value = "'0.2659"
$scratch.Copy()
$ws.Range('D8').PasteSpecial(-4163)
$scratch.Value = "'0.06066"
$scratch.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$scratch.Value = "'0.07169"
$scratch.Copy()
$ws.Range('D11').PasteSpecial(-4163)
$scratch.Value = "'0.6337"
$scratch.Copy()
$ws.Range('D12').PasteSpecial(-4163)
$scratch.Value = "'14.69"
$scratch.Copy()
$ws.Range('D13').PasteSpecial(-4163)
$scratch.Value = "'4.674"
$scratch.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$scratch.Value = "'74.59"
$scratch.Copy()
$ws.Range('D15').PasteSpecial(-4163)
$scratch.Value = "'1.005"
$scratch.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$scratch.Value = "'11.62"
$scratch.Copy()
$ws.Range('D19').PasteSpecial(-4163)
$scratch.Value = "'0.000006703"
$scratch.Copy()
$ws.Range('D20').PasteSpecial(-4163)
$scratch.Value = "'4.498"
$scratch.Copy()
$ws.Range('D21').PasteSpecial(-4163)
$scratch.Value = "'5.360"
$scratch.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$scratch.Value = "'134.04"
$scratch.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$scratch.Value = "'14.99"
$scratch.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$scratch.Value = "'1.399"
$scratch.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$scratch.Value = "'1.741"
$scratch.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$scratch.Value = "'103.60"
$scratch.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$scratch.Value = "'3.859"
$scratch.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$scratch.Value = "'0.08016"
$scratch.Copy()
$ws.Range('D31').PasteSpecial(-4163)
$scratch.Value = "'3.586"
$scratch.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$scratch.Value = "'0.04650"
$scratch.Copy()
$ws.Range('D33').PasteSpecial(-4163)
$scratch.Value = "'2.661"
$scratch.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$scratch.Value = "'0.9698"
$scratch.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$scratch.Value = "'0.5914"
$scratch.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$scratch.Value = "'2.681"
$scratch.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$scratch.Value = "'0.01574"
$scratch.Copy()
$ws.Range('D38').PasteSpecial(-4163)
$scratch.Value = "'0.8465"
$scratch.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$scratch.Value = "'1.005"
$scratch.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$scratch.Value = "'1.897"
$scratch.Copy()
$ws.Range('D41').PasteSpecial(-4163)
$scratch.Value = "'100.11"
$scratch.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$scratch.Value = "'0.3789"
$scratch.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$scratch.Value = "'4.933"
$scratch.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$scratch.Value = "'0.1159"
$scratch.Copy()
$ws.Range('D45').PasteSpecial(-4163)
$scratch.Value = "'6.155"
$scratch.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$scratch.Value = "'54.45"
$scratch.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$scratch.Value = "'0.05218"
$scratch.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$scratch.Value = "'30.00"
$scratch.Copy()
$ws.Range('D49').PasteSpecial(-4163)
$scratch.Value = "'7.485"
$scratch.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$scratch.Value = "'0.3387"
$scratch.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$scratch.EntireColumn.Delete()
$excel.CutCopyMode = 0
